$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "652"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1505687.79"
$ws.Range("D2").Style = "Normal"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "1014"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "3577453.47"
$ws.Range("D4").Style = "Normal"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "644"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2066277.78"
$ws.Range("D6").Style = "Normal"

$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "180"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "457809.00"
$ws.Range("D20").Style = "Normal"

$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "334"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1181741.00"
$ws.Range("D21").Style = "Normal"

$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "159"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "452312.39"
$ws.Range("D22").Style = "Normal"

$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "277"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "718237.45"
$ws.Range("D28").Style = "Normal"

$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "549"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2244182.70"
$ws.Range("D30").Style = "Normal"

$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "383"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1300613.57"
$ws.Range("D32").Style = "Normal"

$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "308"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "732119.71"
$ws.Range("D35").Style = "Normal"

$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "375"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1017567.74"
$ws.Range("D45").Style = "Normal"

$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "612"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2400292.99"
$ws.Range("D47").Style = "Normal"

$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "411"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1393954.16"
$ws.Range("D48").Style = "Normal"

$ws.Range("C59").NumberFormat = "@"
$ws.Range("C59").Value = "49"
$ws.Range("C59").Style = "Normal"
$ws.Range("D59").NumberFormat = "@"
$ws.Range("D59").Value = "113416.69"
$ws.Range("D59").Style = "Normal"

$ws.Range("C60").NumberFormat = "@"
$ws.Range("C60").Value = "89"
$ws.Range("C60").Style = "Normal"
$ws.Range("D60").NumberFormat = "@"
$ws.Range("D60").Value = "271019.00"
$ws.Range("D60").Style = "Normal"

$ws.Range("C61").NumberFormat = "@"
$ws.Range("C61").Value = "23"
$ws.Range("C61").Style = "Normal"
$ws.Range("D61").NumberFormat = "@"
$ws.Range("D61").Value = "74777.00"
$ws.Range("D61").Style = "Normal"

$ws.Range("C74").NumberFormat = "@"
$ws.Range("C74").Value = "386"
$ws.Range("C74").Style = "Normal"
$ws.Range("D74").NumberFormat = "@"
$ws.Range("D74").Value = "977509.70"
$ws.Range("D74").Style = "Normal"

$ws.Range("C76").NumberFormat = "@"
$ws.Range("C76").Value = "924"
$ws.Range("C76").Style = "Normal"
$ws.Range("D76").NumberFormat = "@"
$ws.Range("D76").Value = "3213775.26"
$ws.Range("D76").Style = "Normal"

$ws.Range("C77").NumberFormat = "@"
$ws.Range("C77").Value = "522"
$ws.Range("C77").Style = "Normal"
$ws.Range("D77").NumberFormat = "@"
$ws.Range("D77").Value = "1722125.47"
$ws.Range("D77").Style = "Normal"

$ws.Range("C79").NumberFormat = "@"
$ws.Range("C79").Value = "35"
$ws.Range("C79").Style = "Normal"
$ws.Range("D79").NumberFormat = "@"
$ws.Range("D79").Value = "146180.27"
$ws.Range("D79").Style = "Normal"

$ws.Range("C92").NumberFormat = "@"
$ws.Range("C92").Value = "625"
$ws.Range("C92").Style = "Normal"
$ws.Range("D92").NumberFormat = "@"
$ws.Range("D92").Value = "1525349.94"
$ws.Range("D92").Style = "Normal"

$ws.Range("C94").NumberFormat = "@"
$ws.Range("C94").Value = "1097"
$ws.Range("C94").Style = "Normal"
$ws.Range("D94").NumberFormat = "@"
$ws.Range("D94").Value = "3729249.80"
$ws.Range("D94").Style = "Normal"

$ws.Range("C96").NumberFormat = "@"
$ws.Range("C96").Value = "1016"
$ws.Range("C96").Style = "Normal"
$ws.Range("D96").NumberFormat = "@"
$ws.Range("D96").Value = "3120636.31"
$ws.Range("D96").Style = "Normal"

$ws.Range("C98").NumberFormat = "@"
$ws.Range("C98").Value = "48"
$ws.Range("C98").Style = "Normal"
$ws.Range("D98").NumberFormat = "@"
$ws.Range("D98").Value = "179835.52"
$ws.Range("D98").Style = "Normal"
